$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column D as Text first so numeric-looking price strings
# (e.g. "584.72", "0.999") are stored as literal text, matching the
# original inlineStr cells, instead of being coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.708.87"
$ws.Range("E2").Value = "  +2.47%  "
$ws.Range("D3").Value = "3.396.81"
$ws.Range("E3").Value = "  +1.89%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "584.72"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "180.34"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "0.200"
$ws.Range("E9").Value = "  +9.64%  "
$ws.Range("D10").Value = "0.593"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("D11").Value = "48.39"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").Value = "0.0000287"
$ws.Range("E12").Value = "  +4.48%  "
$ws.Range("D13").Value = "685.02"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("D15").Value = "3.942.14"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "69.680.21"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.398.57"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.121"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "17.72"
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").Value = "17.27"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").Value = "5.35"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").Value = "102.59"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "3.93"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "2.72"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").Value = "33.88"
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("D29").Value = "8.85"
$ws.Range("E29").Value = "  +3.48%  "
$ws.Range("D30").Value = "6.97"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "11.12"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("B32").Value = "dogwifhat"
$ws.Range("C32").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D32").Value = "3.65"
$ws.Range("E32").Value = "  +8.23%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "554.55"
$ws.Range("E33").Value = "  -3.52%  "
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").Value = "58.40"
$ws.Range("E35").Value = "  +2.23%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "3.666.40"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("E38").Value = "  +4.54%  "
$ws.Range("D39").Value = "35.65"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").Value = "0.0₃0736"
$ws.Range("E40").Value = "  +8.67%  "
$ws.Range("D41").Value = "3.30"
$ws.Range("E41").Value = "  +4.13%  "
$ws.Range("E42").Value = "  +2.98%  "
$ws.Range("E43").Value = "  +4.64%  "
$ws.Range("D44").Value = "0.338"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").Value = "3.35"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "1.38"
$ws.Range("E49").Value = "  +3.93%  "
$ws.Range("D50").Value = "129.82"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "2.63"
$ws.Range("E51").Value = "  +0.80%  "
